# Apply updated timetable cell values (adds teacher initials in brackets)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4" = "CS802A[BDu]  /  CS802B[BDu]"; "B4" = "Free Period!"; "C4" = "Free Period!"; "D4" = "Free Period!"; "E4" = "Free Period!"; "F4" = "HU801[KB]  /  "; "G4" = "CS801A[SSK]  /  "
    "A6" = "CS894[MDu]  /  "; "B6" = "CS894[MDu]  /  "; "C6" = "CS894[MDu]  /  "; "D6" = "CS802A[BDu]  /  CS802B[BDu]"; "E6" = "CS891[DC, SMa]  /  "; "F6" = "CS891[DC, SMa]  /  "; "G6" = "CS891[DC, SMa]  /  "
    "A8" = "Free Period!"; "B8" = "CS894[MDu]  /  "; "C8" = "CS894[MDu]  /  "; "D8" = "CS894[MDu]  /  "; "E8" = "Free Period!"; "F8" = "HU801[KB]  /  "; "G8" = "CS801A[SSK]  /  "
    "A10" = "CS802A[BDu]  /  CS802B[BDu]"; "B10" = "CS894[SMa]  /  "; "C10" = "CS894[SMa]  /  "; "D10" = "CS894[SMa]  /  "; "E10" = "CS891[DC, SMa]  /  "; "F10" = "CS891[DC, SMa]  /  "; "G10" = "CS891[DC, SMa]  /  "
    "A12" = "Free Period!"; "B12" = "CS894[GY]  /  "; "C12" = "CS894[GY]  /  "; "D12" = "CS894[GY]  /  "; "E12" = "Free Period!"; "F12" = "Free Period!"; "G12" = "CS801A[SSK]  /  "
    "A16" = "IT802A[KDa]  /  "; "B16" = "IT894[AKS]  /  "; "C16" = "IT894[AKS]  /  "; "D16" = "IT894[AKS]  /  "; "E16" = "Free Period!"; "F16" = "IT801A[RCh]  /  "; "G16" = "Free Period!"
    "A18" = "IT894[AKS]  /  "; "B18" = "IT894[AKS]  /  "; "C18" = "IT894[AKS]  /  "; "D18" = "Free Period!"; "E18" = "IT891[KDa, SU]  /  "; "F18" = "IT891[KDa, SU]  /  "; "G18" = "IT891[KDa, SU]  /  "
    "A20" = "IT801A[RCh]  /  "; "B20" = "Free Period!"; "C20" = "Free Period!"; "D20" = "Free Period!"; "E20" = "IT894[AB]  /  "; "F20" = "IT894[AB]  /  "; "G20" = "IT894[AB]  /  "
    "A22" = "IT802A[KDa]  /  "; "B22" = "Free Period!"; "C22" = "HU801[KB]  /  "; "D22" = "Free Period!"; "E22" = "IT891[KDa, AGh]  /  "; "F22" = "IT891[KDa, AGh]  /  "; "G22" = "IT891[KDa, AGh]  /  "
    "A24" = "IT802A[KDa]  /  "; "B24" = "Free Period!"; "C24" = "HU801[KB]  /  "; "D24" = "IT801A[RCh]  /  "; "E24" = "IT894[RG]  /  "; "F24" = "IT894[RG]  /  "; "G24" = "IT894[RG]  /  "
    "A28" = "HU801[KB]  /  "; "B28" = "Free Period!"; "C28" = "ECE801A[DK]  /  "; "D28" = "Free Period!"; "E28" = "ECE894[TD]  /  "; "F28" = "ECE894[TD]  /  "; "G28" = "ECE894[TD]  /  "
    "A30" = "ECE801A[DK]  /  "; "B30" = "ECE894[RND]  /  "; "C30" = "ECE894[RND]  /  "; "D30" = "ECE894[RND]  /  "; "E30" = "ECE802A[ArD]  /  "; "F30" = "Free Period!"; "G30" = "Free Period!"
    "A32" = "ECE802A[ArD]  /  "; "B32" = "Free Period!"; "C32" = "Free Period!"; "D32" = "Free Period!"; "E32" = "ECE891[SD, AnC]  /  "; "F32" = "ECE891[SD, AnC]  /  "; "G32" = "ECE891[SD, AnC]  /  "
    "A34" = "ECE891[AnC, BC]  /  "; "B34" = "ECE891[AnC, BC]  /  "; "C34" = "ECE891[AnC, BC]  /  "; "D34" = "Free Period!"; "E34" = "ECE894[RND]  /  "; "F34" = "ECE894[RND]  /  "; "G34" = "ECE894[RND]  /  "
    "A36" = "ECE802A[ArD]  /  "; "B36" = "ECE894[PC]  /  "; "C36" = "ECE894[PC]  /  "; "D36" = "ECE894[PC]  /  "; "E36" = "ECE801A[DK]  /  "; "F36" = "HU801[KB]  /  "; "G36" = "Free Period!"
    "A40" = "EE801A[ASG]  /  "; "B40" = "Free Period!"; "C40" = "Free Period!"; "D40" = "HU801[KB]  /  "; "E40" = "EE894[MB]  /  "; "F40" = "EE894[MB]  /  "; "G40" = "EE894[MB]  /  "
    "A42" = "Free Period!"; "B42" = "Free Period!"; "C42" = "EE802A[KR]  /  "; "D42" = "Free Period!"; "E42" = "Free Period!"; "F42" = "EE801A[ASG]  /  "; "G42" = "Free Period!"
    "A44" = "HU801[KB]  /  "; "B44" = "EE894[MB]  /  "; "C44" = "EE894[MB]  /  "; "D44" = "EE894[MB]  /  "; "E44" = "Free Period!"; "F44" = "Free Period!"; "G44" = "EE802A[KR]  /  "
    "A46" = "EE801A[ASG]  /  "; "B46" = "EE891[PG, IB]  /  "; "C46" = "EE891[PG, IB]  /  "; "D46" = "EE891[PG, IB]  /  "; "E46" = "EE894[MB]  /  "; "F46" = "EE894[MB]  /  "; "G46" = "EE894[MB]  /  "
    "A48" = "EE891[PG, IB]  /  "; "B48" = "EE891[PG, IB]  /  "; "C48" = "EE891[PG, IB]  /  "; "D48" = "EE802A[KR]  /  "; "E48" = "EE894[MB]  /  "; "F48" = "EE894[MB]  /  "; "G48" = "EE894[MB]  /  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}